$d = $word.ActiveDocument

# Replace all occurrences of "MiliQ" with "MilliQ" (appears 6 times, same formatting context)
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("MiliQ", $true, $false, $false, $false, $false, $true, 1, $false, "MilliQ", 2)

# Update the revision date from 2022-12-08 to 2022-12-09 in the change log line
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("za: Jolanta Walkusz-Miotk 2022-12-08", $true, $false, $false, $false, $false, $true, 1, $false, "za: Jolanta Walkusz-Miotk 2022-12-09", 2)
